$d = $word.ActiveDocument
$t = $d.Tables(1)

# Header row (row 1): rename columns 3, 4, 5, 6
$t.Cell(1, 3).Range.Text = "Mã Đội Bóng Cũ"
$t.Cell(1, 4).Range.Text = "Tên Cầu Thủ"
$t.Cell(1, 5).Range.Text = "Giá Trị"
$t.Cell(1, 6).Range.Text = "Tên Đội Bóng Cũ"

# Data row (row 2): update values
$t.Cell(2, 1).Range.Text = "9"
$t.Cell(2, 2).Range.Text = "6"
$t.Cell(2, 3).Range.Text = "6"
$t.Cell(2, 5).Range.Text = "1000.0"
$t.Cell(2, 6).Range.Text = "ád"
$t.Cell(2, 7).Range.Text = "27/05/2024"
